# T1753 test data: rename the old "Repeaters" sheet to "RepeatersOld" and
# promote "Repeaters_Updated" to be the new "Repeaters" sheet, then make the
# (new) "Repeaters" sheet the active tab with cell A16 selected.

$wb = $excel.ActiveWorkbook

$sheetOld = $wb.Worksheets.Item(1)   # was "Repeaters"          (sheetId 7)
$sheetNew = $wb.Worksheets.Item(2)   # was "Repeaters_Updated"  (sheetId 11)

$sheetOld.Name = "RepeatersOld"
$sheetNew.Name = "Repeaters"

$sheetNew.Activate()
[void]$sheetNew.Range("A16").Select()
